$wb = $excel.ActiveWorkbook

# --- Sheet: "Block Drops" (sheet1) ---------------------------------------
$blockDrops = $wb.Worksheets.Item("Block Drops")

# Row 91 already has A91 ("Birch Log"); append two more drop items.
$blockDrops.Range("B91").Value = "Weeping Vines (ST)"
$blockDrops.Range("C91").Value = "Snowballs"

# New row 95
$blockDrops.Range("A95").Value = "Obsidian"
$blockDrops.Range("B95").Value = "Oxidized Cut Copper Stairs"
$blockDrops.Range("C95").Value = "Brown Stained Glass Pane"
$blockDrops.Range("D95").Value = "Spruce Pressure Plate"
$blockDrops.Range("E95").Value = "Block of Raw Gold"
$blockDrops.Range("F95").Value = "Cyan Glazed Terracotta"
$blockDrops.Range("G95").Value = "Farm Loot"

# New row 96
$blockDrops.Range("A96").Value = "Red Wool"
$blockDrops.Range("B96").Value = "Nether Loot Chest"

# New row 97
$blockDrops.Range("A97").Value = "Polished Diorite"
$blockDrops.Range("B97").Value = "Crimson Fence"
$blockDrops.Range("C97").Value = "Polished Deepslate Slab"
$blockDrops.Range("D97").Value = "Sea Lantern"
$blockDrops.Range("E97").Value = "Diorite Stairs"
$blockDrops.Range("F97").Value = "Brown Candle"
$blockDrops.Range("G97").Value = "Chiseled Polished Blackstone"
$blockDrops.Range("H97").Value = "Orange Stained Glass"
$blockDrops.Range("I97").Value = "Magenta Glazed Terracotta"
$blockDrops.Range("J97").Value = "Black Glazed Terrcotta"
$blockDrops.Range("K97").Value = "Lily Pad"
$blockDrops.Range("L97").Value = "Lodestone"
$blockDrops.Range("M97").Value = "String"

# --- Sheet: "Mob Drops" (sheet2) -----------------------------------------
$mobDrops = $wb.Worksheets.Item("Mob Drops")
$mobDrops.Range("A15").Select()

# --- Sheet: "Crafting" (sheet3) ------------------------------------------
$crafting = $wb.Worksheets.Item("Crafting")

$crafting.Range("A21").Value = "4 Birch Logs"
$crafting.Range("B21").Value = "4 Polished Diorite"

$crafting.Range("A22").Value = "3 Polished Diorite"
$crafting.Range("B22").Value = "Piston"

$crafting.Range("A23").Value = "4 Snowball"
$crafting.Range("B23").Value = "8 Cyan Concrete Powder"

$crafting.Range("A24").Value = "9 Copper Ingot"
$crafting.Range("B24").Value = "Barrel"

# --- View state -----------------------------------------------------------
# Scroll "Block Drops" so row 87 is at the top, and leave the selection on E63.
$blockDrops.Activate()
$blockDrops.Range("E63").Select()
$excel.ActiveWindow.ScrollRow = 87
$excel.ActiveWindow.ScrollColumn = 1

# "Crafting" ends up the active/selected tab, with C24 selected.
$crafting.Activate()
$crafting.Range("C24").Select()
